$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.103
$ws.Range("E8").Value = 13.419
$ws.Range("A12").Value = -21.992
$ws.Range("E12").Value = 13.143
$ws.Range("E14").Value = 13.072
$ws.Range("E22").Value = 13.107
